# Excel template bug fix
#
# "Bico" sheet: the Obs_relatorio column (H, rows 2-9) was incorrectly
# pre-filled with the "Validado com sucesso!" success message for every
# row regardless of the actual validation outcome - clear it back out.
#
# "Tanque" sheet: the Obs_relatorio column (F, rows 2-5) should report the
# real divergence between the SPED closing value and the report value
# instead of always showing the "Validado com sucesso!" success message.

$wb = $excel.ActiveWorkbook

# --- "Bico" sheet: clear the bogus success message in H2:H9 ---
$bico = $wb.Worksheets.Item("Bico")
for ($row = 2; $row -le 9; $row++) {
    $bico.Cells.Item($row, 8).Value = ""
}

# --- "Tanque" sheet: replace success message with real divergence text ---
$tanque = $wb.Worksheets.Item("Tanque")
$tanque.Cells.Item(2, 6).Value = "Divergência entre o SPED(4397,41) e o relatório(718651,00)!"
$tanque.Cells.Item(3, 6).Value = "Divergência entre o SPED(4397,41) e o relatório(325178,00)!"
$tanque.Cells.Item(4, 6).Value = "Divergência entre o SPED(4397,41) e o relatório(328364,00)!"
$tanque.Cells.Item(5, 6).Value = "Divergência entre o SPED(4397,41) e o relatório(439741,00)!"

Write-Output "Applied Bico H2:H9 clear and Tanque F2:F5 divergence messages"
